$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 20981.428
$ws.Range("J93").Value = 20981.428
$ws.Range("L93").Value = 20981.428
$ws.Range("N93").Value = -25973.428
$ws.Range("H129").Value = 1383.0741
$ws.Range("J129").Value = 1475.92
$ws.Range("L129").Value = 4427.76
$ws.Range("N129").Value = -14427.76
$ws.Range("H132").Value = 36215604
$ws.Range("I132").Value = 52838956
$ws.Range("J132").Value = 1121865.8
$ws.Range("K132").Value = 158516868
$ws.Range("L132").Value = 3365597.4
$ws.Range("M132").Value = -158514338
$ws.Range("N132").Value = -3370657.4
$ws.Range("H137").Value = 901754
$ws.Range("I137").Value = 1987619.9
$ws.Range("J137").Value = 3106.3794
$ws.Range("K137").Value = 5962859.699999999
$ws.Range("L137").Value = 9319.138199999999
$ws.Range("M137").Value = -5960309.699999999
$ws.Range("N137").Value = -14419.1382

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1456.7222
$ws.Range("I61").Value = 1401.3125
$ws.Range("J61").Value = 1900
$ws.Range("K61").Value = 1401.3125
$ws.Range("L61").Value = 1900
$ws.Range("M61").Value = -1189.3125
$ws.Range("N61").Value = -2324
$ws.Range("H74").Value = 411598.38
$ws.Range("I74").Value = 762131.4399999999
$ws.Range("J74").Value = 2643.1667
$ws.Range("K74").Value = 762131.4399999999
$ws.Range("L74").Value = 2643.1667
$ws.Range("M74").Value = -761257.4399999999
$ws.Range("N74").Value = -4391.1667
$ws.Range("H77").Value = 411598.38
$ws.Range("I77").Value = 762131.4399999999
$ws.Range("J77").Value = 2643.1667
$ws.Range("K77").Value = 3810657.2
$ws.Range("L77").Value = 13215.8335
$ws.Range("M77").Value = -3806289.2
$ws.Range("N77").Value = -21951.8335
$ws.Range("H109").Value = 900759.9
$ws.Range("J109").Value = 900759.9
$ws.Range("L109").Value = 900759.9
$ws.Range("N109").Value = -903533.9
$ws.Range("H136").Value = 1456.7222
$ws.Range("I136").Value = 1401.3125
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 4203.9375
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -1653.9375
$ws.Range("N136").Value = -10800
$ws.Range("H137").Value = 48560
$ws.Range("J137").Value = 48560
$ws.Range("L137").Value = 48560
$ws.Range("N137").Value = -58760
$ws.Range("H139").Value = 43070
$ws.Range("J139").Value = 43070
$ws.Range("L139").Value = 43070
$ws.Range("N139").Value = -53350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4205.4287
$ws.Range("I134").Value = 1845.5
$ws.Range("J134").Value = 5657.6924
$ws.Range("K134").Value = 5536.5
$ws.Range("L134").Value = 16973.0772
$ws.Range("M134").Value = -3001.5
$ws.Range("N134").Value = -22043.0772
$ws.Range("H137").Value = 50630
$ws.Range("J137").Value = 50630
$ws.Range("L137").Value = 50630
$ws.Range("N137").Value = -60830
$ws.Range("H138").Value = 39609.03
$ws.Range("J138").Value = 39609.03
$ws.Range("L138").Value = 39609.03
$ws.Range("N138").Value = -49889.03

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 29999
$ws.Range("J11").Value = 29999
$ws.Range("L11").Value = 29999
$ws.Range("N11").Value = -30279
$ws.Range("H22").Value = 602.0625
$ws.Range("I22").Value = 233.16667
$ws.Range("J22").Value = 823.4
$ws.Range("K22").Value = 233.16667
$ws.Range("L22").Value = 823.4
$ws.Range("M22").Value = 116.83333
$ws.Range("N22").Value = -1523.4
$ws.Range("H31").Value = 223663.7
$ws.Range("I31").Value = 501405.75
$ws.Range("J31").Value = 3103.853
$ws.Range("K31").Value = 501405.75
$ws.Range("L31").Value = 3103.853
$ws.Range("M31").Value = -501110.75
$ws.Range("N31").Value = -3693.853
$ws.Range("H34").Value = 223663.7
$ws.Range("I34").Value = 501405.75
$ws.Range("J34").Value = 3103.853
$ws.Range("K34").Value = 501405.75
$ws.Range("L34").Value = 3103.853
$ws.Range("M34").Value = -501203.75
$ws.Range("N34").Value = -3507.853
$ws.Range("H58").Value = 2486.7942
$ws.Range("I58").Value = 1363.8462
$ws.Range("J58").Value = 6136.375
$ws.Range("K58").Value = 1363.8462
$ws.Range("L58").Value = 6136.375
$ws.Range("M58").Value = -1160.8462
$ws.Range("N58").Value = -6542.375
$ws.Range("H107").Value = 1135.5294
$ws.Range("I107").Value = 782.8333
$ws.Range("J107").Value = 1327.909
$ws.Range("K107").Value = 782.8333
$ws.Range("L107").Value = 1327.909
$ws.Range("M107").Value = 1137.1667
$ws.Range("N107").Value = -5167.909
$ws.Range("H132").Value = 4333.5625
$ws.Range("I132").Value = 3665.9285
$ws.Range("J132").Value = 9007
$ws.Range("K132").Value = 10997.7855
$ws.Range("L132").Value = 27021
$ws.Range("M132").Value = -8467.7855
$ws.Range("N132").Value = -32081
$ws.Range("H136").Value = 2486.7942
$ws.Range("I136").Value = 1363.8462
$ws.Range("J136").Value = 6136.375
$ws.Range("K136").Value = 4091.5386
$ws.Range("L136").Value = 18409.125
$ws.Range("M136").Value = -1541.5386
$ws.Range("N136").Value = -23509.125
$ws.Range("H138").Value = 43546.363
$ws.Range("J138").Value = 43546.363
$ws.Range("L138").Value = 43546.363
$ws.Range("N138").Value = -53826.363
$ws.Range("H140").Value = 120705
$ws.Range("J140").Value = 120705
$ws.Range("L140").Value = 120705
$ws.Range("N140").Value = -131065
$ws.Range("H141").Value = 24530
$ws.Range("J141").Value = 24530
$ws.Range("L141").Value = 24530
$ws.Range("N141").Value = -34890

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 944.3857400000001
$ws.Range("I68").Value = 657.6842
$ws.Range("K68").Value = 1973.0526
$ws.Range("M68").Value = -1162.0526
$ws.Range("H69").Value = 1819.5
$ws.Range("I69").Value = 695.5454999999999
$ws.Range("J69").Value = 4292.2
$ws.Range("K69").Value = 2086.6365
$ws.Range("L69").Value = 12876.6
$ws.Range("M69").Value = -1275.6365
$ws.Range("N69").Value = -14498.6
$ws.Range("H71").Value = 944.3857400000001
$ws.Range("I71").Value = 657.6842
$ws.Range("K71").Value = 5919.1578
$ws.Range("M71").Value = -1863.1578
$ws.Range("H72").Value = 1819.5
$ws.Range("I72").Value = 695.5454999999999
$ws.Range("J72").Value = 4292.2
$ws.Range("K72").Value = 6259.9095
$ws.Range("L72").Value = 38629.8
$ws.Range("M72").Value = -2203.9095
$ws.Range("N72").Value = -46741.8
$ws.Range("H113").Value = 733.6957
$ws.Range("I113").Value = 480.5625
$ws.Range("J113").Value = 1312.2858
$ws.Range("K113").Value = 1441.6875
$ws.Range("L113").Value = 3936.8574
$ws.Range("M113").Value = 728.3125
$ws.Range("N113").Value = -8276.857400000001
$ws.Range("H131").Value = 775.49
$ws.Range("J131").Value = 805.10986
$ws.Range("L131").Value = 2415.32958
$ws.Range("N131").Value = -12495.32958
$ws.Range("H132").Value = 4632.5835
$ws.Range("I132").Value = 631.3333
$ws.Range("J132").Value = 5966.3335
$ws.Range("K132").Value = 5681.9997
$ws.Range("L132").Value = 53697.0015
$ws.Range("M132").Value = -3151.9997
$ws.Range("N132").Value = -58757.0015
$ws.Range("H137").Value = 2783.9778
$ws.Range("I137").Value = 1924
$ws.Range("J137").Value = 3471.96
$ws.Range("K137").Value = 5772
$ws.Range("L137").Value = 10415.88
$ws.Range("M137").Value = -672
$ws.Range("N137").Value = -20615.88

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23413.3
$ws.Range("J46").Value = 24348.666
$ws.Range("L46").Value = 24348.666
$ws.Range("N46").Value = -24660.666
$ws.Range("H70").Value = 6068.0835
$ws.Range("I70").Value = 5711.795
$ws.Range("K70").Value = 5711.795
$ws.Range("M70").Value = -5441.795
$ws.Range("H73").Value = 6068.0835
$ws.Range("I73").Value = 5711.795
$ws.Range("K73").Value = 5711.795
$ws.Range("M73").Value = -4775.795
$ws.Range("H132").Value = 3304.342
$ws.Range("I132").Value = 2476.8696
$ws.Range("J132").Value = 4573.1333
$ws.Range("K132").Value = 7430.6088
$ws.Range("L132").Value = 13719.3999
$ws.Range("M132").Value = -4900.6088
$ws.Range("N132").Value = -18779.3999
$ws.Range("H137").Value = 39313.332
$ws.Range("J137").Value = 48970
$ws.Range("L137").Value = 48970
$ws.Range("N137").Value = -59170
$ws.Range("H140").Value = 38363.684
$ws.Range("J140").Value = 38363.684
$ws.Range("L140").Value = 38363.684
$ws.Range("N140").Value = -48723.684

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6800
$ws.Range("I100").Value = 2400
$ws.Range("K100").Value = 2400
$ws.Range("M100").Value = -1859
$ws.Range("H139").Value = 48799.285
$ws.Range("J139").Value = 48799.285
$ws.Range("L139").Value = 48799.285
$ws.Range("N139").Value = -59079.285
$ws.Range("H141").Value = 31928.75
$ws.Range("J141").Value = 31928.75
$ws.Range("L141").Value = 31928.75
$ws.Range("N141").Value = -42288.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 52502.75
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 70002
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 70002
$ws.Range("M2").Value = 107
$ws.Range("N2").Value = -70226
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280
$ws.Range("H139").Value = 42500
$ws.Range("J139").Value = 42500
$ws.Range("L139").Value = 42500
$ws.Range("N139").Value = -52780
$ws.Range("H140").Value = 43931.363
$ws.Range("J140").Value = 46824.5
$ws.Range("L140").Value = 46824.5
$ws.Range("N140").Value = -57184.5
